# Auto-generated Excel COM-interop script to apply the commit diff
# ("chore: update Sheets via scheduled runner") to the workbook.
# Each hunk below corresponds 1:1 to a cell change identified from the
# canonical-OOXML diff, located unambiguously by (sheet, row, old value).

$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(86, 8).Value = 2043  # H86: 1834.6774 -> 2043
$ws.Cells.Item(86, 9).Value = 1974.4615  # I86: 1723.4783 -> 1974.4615
$ws.Cells.Item(86, 11).Value = 1974.4615  # K86: 1723.4783 -> 1974.4615
$ws.Cells.Item(86, 13).Value = -851.4614999999999  # M86: -600.4783 -> -851.4614999999999
$ws.Cells.Item(87, 8).Value = 77000  # H87: 88500 -> 77000
$ws.Cells.Item(87, 10).Value = 0  # J87: 100000 -> 0
$ws.Cells.Item(87, 12).Value = 0  # L87: 100000 -> 0
$ws.Cells.Item(87, 14).ClearContents()  # N87: was -102496
$ws.Cells.Item(89, 8).Value = 2043  # H89: 1834.6774 -> 2043
$ws.Cells.Item(89, 9).Value = 1974.4615  # I89: 1723.4783 -> 1974.4615
$ws.Cells.Item(89, 11).Value = 9872.307499999999  # K89: 8617.3915 -> 9872.307499999999
$ws.Cells.Item(89, 13).Value = -4256.307499999999  # M89: -3001.3915 -> -4256.307499999999
$ws.Cells.Item(90, 8).Value = 77000  # H90: 88500 -> 77000
$ws.Cells.Item(90, 10).Value = 0  # J90: 100000 -> 0
$ws.Cells.Item(90, 12).Value = 0  # L90: 300000 -> 0
$ws.Cells.Item(90, 14).ClearContents()  # N90: was -312480
$ws.Cells.Item(100, 8).Value = 1643.1666  # H100: 1850.5 -> 1643.1666
$ws.Cells.Item(100, 9).Value = 1711.75  # I100: 1919.4 -> 1711.75
$ws.Cells.Item(100, 11).Value = 1711.75  # K100: 1919.4 -> 1711.75
$ws.Cells.Item(100, 13).Value = -1170.75  # M100: -1378.4 -> -1170.75
$ws.Cells.Item(132, 8).Value = 2152.3696  # H132: 2152.587 -> 2152.3696
$ws.Cells.Item(132, 9).Value = 1324.6216  # I132: 1343.6487 -> 1324.6216
$ws.Cells.Item(132, 10).Value = 5555.3335  # J132: 5478.222 -> 5555.3335
$ws.Cells.Item(132, 11).Value = 3973.8648  # K132: 4030.9461 -> 3973.8648
$ws.Cells.Item(132, 12).Value = 16666.0005  # L132: 16434.666 -> 16666.0005
$ws.Cells.Item(132, 13).Value = -1443.8648  # M132: -1500.9461 -> -1443.8648
$ws.Cells.Item(132, 14).Value = -21726.0005  # N132: -21494.666 -> -21726.0005
$ws.Cells.Item(137, 8).Value = 1807.8572  # H137: 1812 -> 1807.8572
$ws.Cells.Item(137, 10).Value = 1918.125  # J137: 1916 -> 1918.125
$ws.Cells.Item(137, 12).Value = 5754.375  # L137: 5748 -> 5754.375
$ws.Cells.Item(137, 14).Value = -10854.375  # N137: -10848 -> -10854.375
$ws.Cells.Item(138, 8).Value = 3376.6  # H138: 3342.12 -> 3376.6
$ws.Cells.Item(138, 9).Value = 1124.6957  # I138: 1094.3043 -> 1124.6957
$ws.Cells.Item(138, 10).Value = 4049.2468  # J138: 4013.5454 -> 4049.2468
$ws.Cells.Item(138, 11).Value = 3374.0871  # K138: 3282.9129 -> 3374.0871
$ws.Cells.Item(138, 12).Value = 12147.7404  # L138: 12040.6362 -> 12147.7404
$ws.Cells.Item(138, 13).Value = 1765.9129  # M138: 1857.0871 -> 1765.9129
$ws.Cells.Item(138, 14).Value = -22427.7404  # N138: -22320.6362 -> -22427.7404
$ws.Cells.Item(141, 8).Value = 685.26  # H141: 662.75 -> 685.26
$ws.Cells.Item(141, 9).Value = 685.26  # I141: 662.75 -> 685.26
$ws.Cells.Item(141, 11).Value = 2055.78  # K141: 1988.25 -> 2055.78
$ws.Cells.Item(141, 13).Value = 3124.22  # M141: 3191.75 -> 3124.22

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 6423.8  # H32: 7278.675 -> 6423.8
$ws.Cells.Item(32, 9).Value = 3840.7222  # I32: 4511.7334 -> 3840.7222
$ws.Cells.Item(32, 10).Value = 16756.111  # J32: 15579.5 -> 16756.111
$ws.Cells.Item(32, 11).Value = 3840.7222  # K32: 4511.7334 -> 3840.7222
$ws.Cells.Item(32, 12).Value = 16756.111  # L32: 15579.5 -> 16756.111
$ws.Cells.Item(32, 13).Value = -3553.7222  # M32: -4224.7334 -> -3553.7222
$ws.Cells.Item(32, 14).Value = -17330.111  # N32: -16153.5 -> -17330.111
$ws.Cells.Item(44, 8).Value = 29833  # H44: 34750 -> 29833
$ws.Cells.Item(44, 10).Value = 34999.5  # J44: 50000 -> 34999.5
$ws.Cells.Item(44, 12).Value = 34999.5  # L44: 50000 -> 34999.5
$ws.Cells.Item(44, 14).Value = -35975.5  # N44: -50976 -> -35975.5
$ws.Cells.Item(55, 8).Value = 10249.5  # H55: 6999.6665 -> 10249.5
$ws.Cells.Item(55, 10).Value = 19999  # J55: 0 -> 19999
$ws.Cells.Item(55, 12).Value = 19999  # L55: 0 -> 19999
$ws.Cells.Item(55, 14).Value = -20629  # N55: None -> -20629
$ws.Cells.Item(61, 8).Value = 3263.7073  # H61: 2973.1064 -> 3263.7073
$ws.Cells.Item(61, 9).Value = 3202.3794  # I61: 2777.9722 -> 3202.3794
$ws.Cells.Item(61, 10).Value = 3411.9167  # J61: 3611.7273 -> 3411.9167
$ws.Cells.Item(61, 11).Value = 3202.3794  # K61: 2777.9722 -> 3202.3794
$ws.Cells.Item(61, 12).Value = 3411.9167  # L61: 3611.7273 -> 3411.9167
$ws.Cells.Item(61, 13).Value = -2990.3794  # M61: -2565.9722 -> -2990.3794
$ws.Cells.Item(61, 14).Value = -3835.9167  # N61: -4035.7273 -> -3835.9167
$ws.Cells.Item(88, 8).Value = 1847.5  # H88: 1126 -> 1847.5
$ws.Cells.Item(88, 9).Value = 1006  # I88: 1002.5 -> 1006
$ws.Cells.Item(88, 10).Value = 2015.8  # J88: 1249.5 -> 2015.8
$ws.Cells.Item(88, 11).Value = 1006  # K88: 1002.5 -> 1006
$ws.Cells.Item(88, 12).Value = 2015.8  # L88: 1249.5 -> 2015.8
$ws.Cells.Item(88, 13).Value = -600  # M88: -596.5 -> -600
$ws.Cells.Item(88, 14).Value = -2827.8  # N88: -2061.5 -> -2827.8
$ws.Cells.Item(91, 8).Value = 1847.5  # H91: 1126 -> 1847.5
$ws.Cells.Item(91, 9).Value = 1006  # I91: 1002.5 -> 1006
$ws.Cells.Item(91, 10).Value = 2015.8  # J91: 1249.5 -> 2015.8
$ws.Cells.Item(91, 11).Value = 1006  # K91: 1002.5 -> 1006
$ws.Cells.Item(91, 12).Value = 2015.8  # L91: 1249.5 -> 2015.8
$ws.Cells.Item(91, 13).Value = 398  # M91: 401.5 -> 398
$ws.Cells.Item(91, 14).Value = -4823.8  # N91: -4057.5 -> -4823.8
$ws.Cells.Item(132, 8).Value = 2378.0625  # H132: 2342.3062 -> 2378.0625
$ws.Cells.Item(132, 9).Value = 2018.9459  # I132: 2019.2703 -> 2018.9459
$ws.Cells.Item(132, 10).Value = 3586  # J132: 3338.3333 -> 3586
$ws.Cells.Item(132, 11).Value = 6056.8377  # K132: 6057.810899999999 -> 6056.8377
$ws.Cells.Item(132, 12).Value = 10758  # L132: 10014.9999 -> 10758
$ws.Cells.Item(132, 13).Value = -3526.8377  # M132: -3527.810899999999 -> -3526.8377
$ws.Cells.Item(132, 14).Value = -15818  # N132: -15074.9999 -> -15818
$ws.Cells.Item(136, 8).Value = 3263.7073  # H136: 2973.1064 -> 3263.7073
$ws.Cells.Item(136, 9).Value = 3202.3794  # I136: 2777.9722 -> 3202.3794
$ws.Cells.Item(136, 10).Value = 3411.9167  # J136: 3611.7273 -> 3411.9167
$ws.Cells.Item(136, 11).Value = 9607.138199999999  # K136: 8333.9166 -> 9607.138199999999
$ws.Cells.Item(136, 12).Value = 10235.7501  # L136: 10835.1819 -> 10235.7501
$ws.Cells.Item(136, 13).Value = -7057.138199999999  # M136: -5783.9166 -> -7057.138199999999
$ws.Cells.Item(136, 14).Value = -15335.7501  # N136: -15935.1819 -> -15335.7501

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 2000  # H86: 1611.1111 -> 2000
$ws.Cells.Item(86, 9).Value = 2000  # I86: 1611.1111 -> 2000
$ws.Cells.Item(86, 11).Value = 2000  # K86: 1611.1111 -> 2000
$ws.Cells.Item(86, 13).Value = -877  # M86: -488.1111000000001 -> -877
$ws.Cells.Item(89, 8).Value = 2000  # H89: 1611.1111 -> 2000
$ws.Cells.Item(89, 9).Value = 2000  # I89: 1611.1111 -> 2000
$ws.Cells.Item(89, 11).Value = 10000  # K89: 8055.5555 -> 10000
$ws.Cells.Item(89, 13).Value = -4384  # M89: -2439.5555 -> -4384
$ws.Cells.Item(99, 8).Value = 3612.3333  # H99: 3484.842 -> 3612.3333
$ws.Cells.Item(99, 10).Value = 5172.3  # J99: 4810.273 -> 5172.3
$ws.Cells.Item(99, 12).Value = 5172.3  # L99: 4810.273 -> 5172.3
$ws.Cells.Item(99, 14).Value = -8168.3  # N99: -7806.273 -> -8168.3
$ws.Cells.Item(107, 8).Value = 8161.067  # H107: 7691.875 -> 8161.067
$ws.Cells.Item(107, 9).Value = 1793.1  # I107: 1689.5454 -> 1793.1
$ws.Cells.Item(107, 11).Value = 1793.1  # K107: 1689.5454 -> 1793.1
$ws.Cells.Item(107, 13).Value = 126.9000000000001  # M107: 230.4546 -> 126.9000000000001
$ws.Cells.Item(134, 8).Value = 1869.9811  # H134: 1732.0333 -> 1869.9811
$ws.Cells.Item(134, 9).Value = 1872.32  # I134: 1745.3572 -> 1872.32
$ws.Cells.Item(134, 10).Value = 1831  # J134: 1545.5 -> 1831
$ws.Cells.Item(134, 11).Value = 5616.96  # K134: 5236.071599999999 -> 5616.96
$ws.Cells.Item(134, 12).Value = 5493  # L134: 4636.5 -> 5493
$ws.Cells.Item(134, 13).Value = -3081.96  # M134: -2701.071599999999 -> -3081.96
$ws.Cells.Item(134, 14).Value = -10563  # N134: -9706.5 -> -10563

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 23831.06  # H31: 23384.44 -> 23831.06
$ws.Cells.Item(31, 9).Value = 37249.207  # I31: 36057.566 -> 37249.207
$ws.Cells.Item(31, 11).Value = 37249.207  # K31: 36057.566 -> 37249.207
$ws.Cells.Item(31, 13).Value = -36954.207  # M31: -35762.566 -> -36954.207
$ws.Cells.Item(34, 8).Value = 23831.06  # H34: 23384.44 -> 23831.06
$ws.Cells.Item(34, 9).Value = 37249.207  # I34: 36057.566 -> 37249.207
$ws.Cells.Item(34, 11).Value = 37249.207  # K34: 36057.566 -> 37249.207
$ws.Cells.Item(34, 13).Value = -37047.207  # M34: -35855.566 -> -37047.207

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(2, 8).Value = 38529.31  # H2: 41732.082 -> 38529.31
$ws.Cells.Item(2, 9).Value = 62558.875  # I2: 71491.14 -> 62558.875
$ws.Cells.Item(2, 10).Value = 82  # J2: 69.40000000000001 -> 82
$ws.Cells.Item(2, 11).Value = 375353.25  # K2: 428946.84 -> 375353.25
$ws.Cells.Item(2, 12).Value = 492  # L2: 416.4 -> 492
$ws.Cells.Item(2, 13).Value = -375240.25  # M2: -428833.84 -> -375240.25
$ws.Cells.Item(2, 14).Value = -718  # N2: -642.4000000000001 -> -718

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 2467.6365  # H97: 2649.4 -> 2467.6365
$ws.Cells.Item(97, 9).Value = 2424.5  # I97: 2678 -> 2424.5
$ws.Cells.Item(97, 11).Value = 2424.5  # K97: 2678 -> 2424.5
$ws.Cells.Item(97, 13).Value = -1928.5  # M97: -2182 -> -1928.5
$ws.Cells.Item(102, 8).Value = 62501204  # H102: 66667884 -> 62501204
$ws.Cells.Item(102, 9).Value = 1275.2858  # I102: 1296.4615 -> 1275.2858
$ws.Cells.Item(102, 11).Value = 1275.2858  # K102: 1296.4615 -> 1275.2858
$ws.Cells.Item(102, 13).Value = 346.7141999999999  # M102: 325.5385000000001 -> 346.7141999999999
$ws.Cells.Item(107, 8).Value = 1515.6316  # H107: 1516.9474 -> 1515.6316
$ws.Cells.Item(107, 9).Value = 1718.4445  # I107: 1619.1 -> 1718.4445
$ws.Cells.Item(107, 10).Value = 1333.1  # J107: 1403.4445 -> 1333.1
$ws.Cells.Item(107, 11).Value = 1718.4445  # K107: 1619.1 -> 1718.4445
$ws.Cells.Item(107, 12).Value = 1333.1  # L107: 1403.4445 -> 1333.1
$ws.Cells.Item(107, 13).Value = 201.5554999999999  # M107: 300.9000000000001 -> 201.5554999999999
$ws.Cells.Item(107, 14).Value = -5173.1  # N107: -5243.4445 -> -5173.1
$ws.Cells.Item(132, 8).Value = 3367.1  # H132: 3247.2188 -> 3367.1
$ws.Cells.Item(132, 9).Value = 3410.7856  # I132: 3280 -> 3410.7856
$ws.Cells.Item(132, 11).Value = 10232.3568  # K132: 9840 -> 10232.3568
$ws.Cells.Item(132, 13).Value = -7702.356800000001  # M132: -7310 -> -7702.356800000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 4646.7417  # H40: 4636.355 -> 4646.7417
$ws.Cells.Item(40, 9).Value = 4114.5415  # I40: 4097.08 -> 4114.5415
$ws.Cells.Item(40, 10).Value = 6471.4287  # J40: 6883.3335 -> 6471.4287
$ws.Cells.Item(40, 11).Value = 4114.5415  # K40: 4097.08 -> 4114.5415
$ws.Cells.Item(40, 12).Value = 6471.4287  # L40: 6883.3335 -> 6471.4287
$ws.Cells.Item(40, 13).Value = -3978.5415  # M40: -3961.08 -> -3978.5415
$ws.Cells.Item(40, 14).Value = -6743.4287  # N40: -7155.3335 -> -6743.4287
$ws.Cells.Item(46, 8).Value = 1080.625  # H46: 1160.5555 -> 1080.625
$ws.Cells.Item(46, 9).Value = 1092.2858  # I46: 1180.75 -> 1092.2858
$ws.Cells.Item(46, 11).Value = 1092.2858  # K46: 1180.75 -> 1092.2858
$ws.Cells.Item(46, 13).Value = -904.2858000000001  # M46: -992.75 -> -904.2858000000001
$ws.Cells.Item(80, 8).Value = 30316.666  # H80: 31190 -> 30316.666
$ws.Cells.Item(80, 10).Value = 30316.666  # J80: 31190 -> 30316.666
$ws.Cells.Item(80, 12).Value = 30316.666  # L80: 31190 -> 30316.666
$ws.Cells.Item(80, 14).Value = -32562.666  # N80: -33436 -> -32562.666
$ws.Cells.Item(83, 8).Value = 30316.666  # H83: 31190 -> 30316.666
$ws.Cells.Item(83, 10).Value = 30316.666  # J83: 31190 -> 30316.666
$ws.Cells.Item(83, 12).Value = 90949.99800000001  # L83: 93570 -> 90949.99800000001
$ws.Cells.Item(83, 14).Value = -102181.998  # N83: -104802 -> -102181.998
$ws.Cells.Item(132, 8).Value = 4149.5186  # H132: 4476.5415 -> 4149.5186
$ws.Cells.Item(132, 9).Value = 3936.5454  # I132: 4205.2 -> 3936.5454
$ws.Cells.Item(132, 10).Value = 5086.6  # J132: 5833.25 -> 5086.6
$ws.Cells.Item(132, 11).Value = 11809.6362  # K132: 12615.6 -> 11809.6362
$ws.Cells.Item(132, 12).Value = 15259.8  # L132: 17499.75 -> 15259.8
$ws.Cells.Item(132, 13).Value = -9279.636200000001  # M132: -10085.6 -> -9279.636200000001
$ws.Cells.Item(132, 14).Value = -20319.8  # N132: -22559.75 -> -20319.8
$ws.Cells.Item(136, 8).Value = 3011.3784  # H136: 2606.2222 -> 3011.3784
$ws.Cells.Item(136, 9).Value = 2783.5862  # I136: 2327.7437 -> 2783.5862
$ws.Cells.Item(136, 10).Value = 3837.125  # J136: 4416.3335 -> 3837.125
$ws.Cells.Item(136, 11).Value = 8350.758600000001  # K136: 6983.2311 -> 8350.758600000001
$ws.Cells.Item(136, 12).Value = 11511.375  # L136: 13249.0005 -> 11511.375
$ws.Cells.Item(136, 13).Value = -5800.758600000001  # M136: -4433.2311 -> -5800.758600000001
$ws.Cells.Item(136, 14).Value = -16611.375  # N136: -18349.0005 -> -16611.375

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 11035.286  # H81: 10492.934 -> 11035.286
$ws.Cells.Item(81, 10).Value = 4766.3335  # J81: 4579.7 -> 4766.3335
$ws.Cells.Item(81, 12).Value = 9532.666999999999  # L81: 9159.4 -> 9532.666999999999
$ws.Cells.Item(81, 14).Value = -11654.667  # N81: -11281.4 -> -11654.667
$ws.Cells.Item(84, 8).Value = 11035.286  # H84: 10492.934 -> 11035.286
$ws.Cells.Item(84, 10).Value = 4766.3335  # J84: 4579.7 -> 4766.3335
$ws.Cells.Item(84, 12).Value = 47663.335  # L84: 45797 -> 47663.335
$ws.Cells.Item(84, 14).Value = -58271.335  # N84: -56405 -> -58271.335
$ws.Cells.Item(107, 8).Value = 803.6129  # H107: 819.2 -> 803.6129
$ws.Cells.Item(107, 9).Value = 815.8461  # I107: 835.04 -> 815.8461
$ws.Cells.Item(107, 11).Value = 2447.5383  # K107: 2505.12 -> 2447.5383
$ws.Cells.Item(107, 13).Value = -527.5383000000002  # M107: -585.1199999999999 -> -527.5383000000002
$ws.Cells.Item(136, 8).Value = 1821.06  # H136: 1789.3726 -> 1821.06
$ws.Cells.Item(136, 9).Value = 1750.0698  # I136: 1714.9546 -> 1750.0698
$ws.Cells.Item(136, 11).Value = 5250.2094  # K136: 5144.8638 -> 5250.2094
$ws.Cells.Item(136, 13).Value = -2700.2094  # M136: -2594.8638 -> -2700.2094

Write-Host "Applied 198 cell updates across 8 sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)."
